$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Belgian Malinois']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new first column (old A:D -> B:E), and two new rows after the
#    current data row (old row 2 stays row 2; new rows 3 and 4 are added for
#    the Samples and Files tabs).
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()
$ws.Rows("3:4").Insert()

# ---------------------------------------------------------------------------
# 2. New "TabName" label column.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# ---------------------------------------------------------------------------
# 3. Populate the Samples and Files query rows, then replace the (shifted)
#    Cases query in B2 with the new, longer query (written last so the
#    shared-string table ends up in the same append order as the source
#    workbook).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery
$ws.Range("B2").Value = $casesQuery

$ws.Range("B2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("B4").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Fill in the StatQuery / dbExcel filename / WebExcel filename columns for
#    the two new rows (copied from row 2, which already holds the right
#    values after the column insert).
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("C4").Value = $ws.Range("C2").Value()
$ws.Range("C3").WrapText = $true
$ws.Range("C4").WrapText = $true

$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("D4").Value = $ws.Range("D2").Value()

$ws.Range("E3").Value = $ws.Range("E2").Value()
$ws.Range("E4").Value = $ws.Range("E2").Value()

# ---------------------------------------------------------------------------
# 5. Column A is brand new, so it needs an explicit width; B:E keep the exact
#    widths they already carried through the column insert.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10

# ---------------------------------------------------------------------------
# 6. Row heights (to fit the long, wrapped query text).
# ---------------------------------------------------------------------------
$ws.Rows("2").RowHeight = 275.5
$ws.Rows("3").RowHeight = 232
$ws.Rows("4").RowHeight = 246.5

# ---------------------------------------------------------------------------
# 7. View: zoom out and move the selection.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 40
$ws.Range("H2").Select()
